# backend and frontend update on 18-12
# Refresh of currentinventory_1_uk_december2025.xlsx: inventory counters and
# their "synced_at" timestamps (column I) move from the 10-Dec sync to the
# 18-Dec sync for every SKU row. Updates total_quantity (D), available_quantity
# (F), reserved_quantity (G), fulfillable_quantity (H), synced_at (I),
# "Current Month Units Sold" (J), "Others" (M) and "Inventory at the end of
# the month" (N) per row, plus the Total row (16) rollups.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    # Row 2  - BV-6X5T-6CY1 / Classic
    "D2" = 4906;  "F2" = 4862;  "H2" = 4862;  "I2" = 46009.22343447226;  "J2" = 4941;  "M2" = 141;  "N2" = 4862;

    # Row 3  - SEMNIWPF / Passion Fruit
    "D3" = 893;   "F3" = 862;   "G3" = 60;    "H3" = 862;   "I3" = 46009.22343448379;  "J3" = 845;   "M3" = 31;   "N3" = 862;

    # Row 4  - SEMNIWME / Menthol
    "D4" = 579;   "F4" = 241;   "G4" = 12;    "H4" = 241;   "I4" = 46009.22343448379;  "J4" = 243;                "N4" = 241;

    # Row 5  - SEWMNIW / Women
    "D5" = 59;    "F5" = 59;    "G5" = 0;     "H5" = 59;    "I5" = 46009.22343448379;  "J5" = 62;                 "N5" = 59;

    # Row 6  - SEWIPES / Intimate Wipes
    "D6" = 294;   "F6" = 273;   "G6" = 10;    "H6" = 273;   "I6" = 46009.22343448379;  "J6" = 265;   "M6" = 16;   "N6" = 273;

    # Row 7  - SEMNIWRF / Refill Pack
    "D7" = 894;   "F7" = 887;   "G7" = 14;    "H7" = 887;   "I7" = 46009.22343448379;  "J7" = 915;   "M7" = 24;   "N7" = 887;

    # Row 8  - SEFMTM / Turmeric
    "D8" = 332;   "F8" = 319;                 "H8" = 319;   "I8" = 46009.22343447226;  "J8" = 320;   "M8" = 5;    "N8" = 319;

    # Row 9  - SECFSH / Shampoo (only the sync timestamp refreshed)
    "I9" = 46009.22343447226;

    # Row 10 - 2Y-LVEI-G7L5 / Classic +Classic
    "D10" = 971;  "F10" = 964;  "G10" = 14;   "H10" = 964;  "I10" = 46009.22343447226; "J10" = 984;               "N10" = 964;

    # Row 11 - SEMNIWCPF / Classic+ Passion Fruit
    "D11" = 1088; "F11" = 1082; "G11" = 12;   "H11" = 1082; "I11" = 46009.22343448379; "J11" = 1091; "M11" = 16;  "N11" = 1082;

    # Row 12 - SEMNIWCM / Classic + Menthol (only the sync timestamp refreshed)
    "I12" = 46009.22343448379;

    # Row 13 - SEIWHCIWWI / Classic + Wipes (only the sync timestamp refreshed)
    "I13" = 46009.22343447226;

    # Row 15 - SEIWHCMEWI / Wipes + Menthol
    "D15" = 142;  "F15" = 142;                "H15" = 142;  "I15" = 46009.22343447226; "J15" = 146;               "N15" = 142;

    # Row 16 - Total row
    "D16" = 10772; "F16" = 10305; "G16" = 210; "H16" = 10305;              "J16" = 10425; "M16" = 253; "N16" = 10305;
}

foreach ($ref in $updates.Keys) {
    $ws.Range($ref).Value = $updates[$ref]
}
